$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Row 21 corresponds to the "CoordinatorDashboard" screen.
# Fill in the "Vervolgschermen" (F) and "Rol" (G) columns.
$ws.Range("F21").Value = "NewQuiz, New Question"
$ws.Range("G21").Value = "Coördinator"

$ws.Range("F21").Select()
